$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 31 with same formatting as row 30 ---
$src = $ws.Range("A30:L30")
$dst = $ws.Range("A31:L31")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Fill in the new record's values ---
$ws.Cells.Item(31, 1).Value = 44016          # A31 DATETIME (7/4/2020)
$ws.Cells.Item(31, 2).Value = 98.7           # B31 WAIST
$ws.Cells.Item(31, 3).Value = 105            # C31 HIP
$ws.Cells.Item(31, 4).Value = 0.94           # D31 WHR
$ws.Cells.Item(31, 5).Value = "SAME"         # E31 WHR_IMPROVEMENT
$ws.Cells.Item(31, 6).Value = 0              # F31 WHR_CHANGE
$ws.Cells.Item(31, 7).Value = 81.9           # G31 WEIGHT
$ws.Cells.Item(31, 8).Value = 0.58           # H31 WHTR
$ws.Cells.Item(31, 9).Value = "SAME"         # I31 WHTR_IMPROVEMENT
$ws.Cells.Item(31, 10).Value = 0             # J31 WHTR_CHANGE
$ws.Cells.Item(31, 11).Value = 28.4          # K31 BMI
$ws.Cells.Item(31, 12).Value = "OVERWEIGHT"  # L31 OBESITY

# --- Restore view state: frozen header row, scrolled down, K31 selected ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$excel.ActiveWindow.ScrollRow = 12
$ws.Range("K31").Select()
